$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (RandomForest)
$ws.Range("B2").Value = 111
$ws.Range("C2").Value = 104
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 19604.28
$ws.Range("F2").Value = 298
$ws.Range("G2").Value = 1282.961634261073
$ws.Range("H2").Value = 93.69
$ws.Range("I2").Value = 4.29
$ws.Range("J2").Value = 887.96
$ws.Range("K2").Value = 177.04
$ws.Range("L2").Value = -55.5
$ws.Range("M2").Value = 199.25
$ws.Range("N2").Value = 1.8
$ws.Range("O2").Value = 18023.32
$ws.Range("P2").Value = 12.4

# Row 3 (XGBoost)
$ws.Range("B3").Value = 111
$ws.Range("C3").Value = 106
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 19870.78
$ws.Range("F3").Value = 221
$ws.Range("G3").Value = 1289.975419289334
$ws.Range("H3").Value = 95.5
$ws.Range("I3").Value = 4.27
$ws.Range("J3").Value = 887.96
$ws.Range("K3").Value = 175.86
$ws.Range("L3").Value = -56.31
$ws.Range("M3").Value = 204
$ws.Range("N3").Value = 1.84
$ws.Range("O3").Value = 18359.8
$ws.Range("P3").Value = 13.15

# Row 4 (Logistic Regression)
$ws.Range("B4").Value = 107
$ws.Range("C4").Value = 105
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 18890.28
$ws.Range("F4").Value = 91
$ws.Range("G4").Value = 1240.120655435238
$ws.Range("H4").Value = 98.13
$ws.Range("I4").Value = 4.12
$ws.Range("J4").Value = 394.09
$ws.Range("K4").Value = 168.37
$ws.Range("L4").Value = -59.91
$ws.Range("M4").Value = 199.25
$ws.Range("N4").Value = 1.86
$ws.Range("O4").Value = 17559.16
$ws.Range("P4").Value = 14.19

# Row 5 (Voting Classifier)
$ws.Range("B5").Value = 110
$ws.Range("C5").Value = 105
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 19690.28
$ws.Range("F5").Value = 221
$ws.Range("G5").Value = 1276.44671611948
$ws.Range("H5").Value = 95.45
$ws.Range("I5").Value = 4.28
$ws.Range("J5").Value = 887.96
$ws.Range("K5").Value = 175.95
$ws.Range("L5").Value = -56.31
$ws.Range("M5").Value = 202.25
$ws.Range("N5").Value = 1.84
$ws.Range("O5").Value = 18192.83
$ws.Range("P5").Value = 13.15

# Row 6 (Stacking Classifier)
$ws.Range("B6").Value = 108
$ws.Range("C6").Value = 105
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 18890.28
$ws.Range("F6").Value = 121
$ws.Range("G6").Value = 1255.098155435238
$ws.Range("H6").Value = 97.22
$ws.Range("I6").Value = 4.12
$ws.Range("J6").Value = 394.09
$ws.Range("K6").Value = 168.37
$ws.Range("L6").Value = -54.93
$ws.Range("M6").Value = 200.75
$ws.Range("N6").Value = 1.86
$ws.Range("O6").Value = 17514.18
$ws.Range("P6").Value = 13.73
